# Update "想去人数" (want-to-go count) figures on the two worksheets that
# list exhibition events: "展览" and "全部类型".
# F2: 632 -> 633
# F3: 478 -> 480
# F8: 1327 -> 1338
# F9: 3968 -> 3981

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 633
    $ws.Range("F3").Value = 480
    $ws.Range("F8").Value = 1338
    $ws.Range("F9").Value = 3981
}
